# Auto-applies the 19:54:54 scrape refresh for Linea 141 (LP1912 / LP1912-215 / 6203-6173).
# Mirrors the upstream scraper re-run: header timestamps bump, a handful of
# existing LP1912 rows get refreshed Hora_Scrap/Linea/Minutos values, and four
# brand-new arrivals are appended at the bottom of the LP1912 table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Header banner: new scrape timestamp / row count -----------------------
$ws1.Range("A2").Value = 'Última actualización: 19:54:54'
$ws1.Range("A3").Value = 'Total filas: 469'
$ws2.Range("A2").Value = 'Última actualización: 19:54:54'
$ws3.Range("A2").Value = 'Última actualización: 19:54:54'

# --- LP1912 data rows refreshed by the new scrape ---------------------------
# Each entry only lists the columns that actually change for that row
# (Hora_Scrap / Linea / Minutos -- Hora_Llegada stays put unless noted).
$cellEdits = @(
    @{Row=48; A='06:33:46'; C='23_HERNANDEZ'; D=87},
    @{Row=49; A='06:16:15'; C='11_ETCHEVERRY'; D=104},
    @{Row=107; A='08:11:27'; C='215C_EL PATO'; D=112},
    @{Row=108; A='09:21:49'; C='23_HERNANDEZ'; D=42},
    @{Row=149; A='10:04:17'; C='23_HERNANDEZ'; D=67},
    @{Row=150; A='10:36:18'; C='15_ABASTO'; D=35},
    @{Row=151; A='11:11:31'; C='16_SANTA ANA'; D=0},
    @{Row=164; A='10:04:17'; C='10_OLMOS'; D=90},
    @{Row=165; A='11:34:25'; C='16_SANTA ANA'; D=0},
    @{Row=166; C='23_HERNANDEZ'},
    @{Row=194; A='10:36:18'; C='11_ETCHEVERRY'; D=114},
    @{Row=195; A='11:53:59'; C='16_P MOR-SANTA ANA'; D=37},
    @{Row=224; C='215_ALUAR'},
    @{Row=225; C='16_SANTA ANA'},
    @{Row=238; A='11:53:59'; C='16_P MOR-SANTA ANA'; D=97},
    @{Row=239; A='11:47:13'; C='10_OLMOS'; D=103},
    @{Row=240; A='11:34:25'; C='215A_EL PATO'; D=116},
    @{Row=293; A='13:39:24'; C='14_ABASTO'; D=94},
    @{Row=294; A='14:31:57'; C='10_OLMOS'; D=42},
    @{Row=381; A='17:13:12'; C='23_HERNANDEZ'; D=59},
    @{Row=382; A='18:10:23'; C='10_OLMOS'; D=2},
    @{Row=392; A='17:54:41'; C='23_HERNANDEZ'; D=39},
    @{Row=393; A='17:34:55'; C='14X44_ABASTO'; D=59},
    @{Row=394; A='16:36:34'; C='14X44_ABASTO'; D=118},
    @{Row=395; A='17:47:22'; C='23_HERNANDEZ'; D=47},
    @{Row=406; C='16_P MOR-SANTA ANA'},
    @{Row=407; C='15_ABASTO'},
    @{Row=444; A='19:54:54'; B='19:54'; C='16_SANTA ANA'; D=0},
    @{Row=445; A='18:30:56'; B='19:58'; D=88},
    @{Row=446; A='18:10:23'; B='19:59'; C='14X44_ABASTO'; D=109},
    @{Row=447; A='18:30:56'; B='20:00'; C='215C_EL PATO'; D=90},
    @{Row=448; A='19:54:54'; B='20:00'; C='16_SANTA ANA'; D=6},
    @{Row=449; B='20:01'; C='16_SANTA ANA'; D=26},
    @{Row=450; A='18:10:23'; B='20:01'; C='215C_EL PATO'; D=111},
    @{Row=451; A='19:35:19'; B='20:03'; C='23_HERNANDEZ'; D=28},
    @{Row=452; A='19:47:42'; B='20:10'; C='23_HERNANDEZ'; D=23},
    @{Row=453; A='18:30:56'; B='20:13'; C='11_ETCHEVERRY'; D=103},
    @{Row=454; A='18:44:14'; B='20:14'; C='11_ETCHEVERRY'; D=90},
    @{Row=455; A='19:47:42'; B='20:15'; C='23_HERNANDEZ'; D=28},
    @{Row=456; A='19:54:54'; B='20:21'; C='23_HERNANDEZ'; D=27},
    @{Row=457; A='18:30:56'; B='20:25'; C='15_ABASTO'; D=115},
    @{Row=458; A='18:44:14'; B='20:26'; C='15_ABASTO'; D=102},
    @{Row=459; A='18:30:56'; B='20:28'; C='10_OLMOS'; D=118},
    @{Row=460; A='18:44:14'; B='20:29'; C='10_OLMOS'; D=105},
    @{Row=461; A='19:11:56'; B='20:43'; C='215B_EL PATO'; D=92},
    @{Row=462; B='20:44'; C='215B_EL PATO'; D=112},
    @{Row=463; B='20:44'; C='17X38_ROMERO'; D=93},
    @{Row=464; A='18:52:19'; B='20:45'; C='17X38_ROMERO'; D=113},
    @{Row=465; A='18:52:19'; B='20:49'; C='11_ETCHEVERRY'; D=117},
    @{Row=466; B='20:51'; C='11_ETCHEVERRY'; D=100},
    @{Row=467; A='19:47:42'; B='20:55'; C='27_EL RETIRO'; D=68},
    @{Row=468; A='19:11:56'; B='20:56'; C='27_EL RETIRO'; D=105},
    @{Row=469; A='19:11:56'; B='21:01'; C='215A_EL PATO'; D=110},
    @{Row=470; B='21:02'; C='215A_EL PATO'; D=87}
)

foreach ($edit in $cellEdits) {
    if ($edit.ContainsKey('A')) { $ws1.Cells.Item($edit.Row, 1).Value = $edit.A }
    if ($edit.ContainsKey('B')) { $ws1.Cells.Item($edit.Row, 2).Value = $edit.B }
    if ($edit.ContainsKey('C')) { $ws1.Cells.Item($edit.Row, 3).Value = $edit.C }
    if ($edit.ContainsKey('D')) { $ws1.Cells.Item($edit.Row, 4).Value = $edit.D }
}

# --- Four brand-new LP1912 arrivals appended at rows 471-474 ----------------
$newRows = @(
    @{Row=471; A='19:47:42'; B='21:09'; C='11_ETCHEVERRY'; D=82; E='LP1912'},
    @{Row=472; A='19:47:42'; B='21:23'; C='10_OLMOS'; D=96; E='LP1912'},
    @{Row=473; A='19:35:19'; B='21:24'; C='10_OLMOS'; D=109; E='LP1912'},
    @{Row=474; A='19:54:54'; B='21:48'; C='11_ETCHEVERRY'; D=114; E='LP1912'}
)

foreach ($row in $newRows) {
    $ws1.Cells.Item($row.Row, 1).Value = $row.A
    $ws1.Cells.Item($row.Row, 2).Value = $row.B
    $ws1.Cells.Item($row.Row, 3).Value = $row.C
    $ws1.Cells.Item($row.Row, 4).Value = $row.D
    $ws1.Cells.Item($row.Row, 5).Value = $row.E
}
